$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C3 currently uses the red font style (fontId=1). Change its font color back to
# automatic/default (black) while keeping the same Calibri font - this creates a
# new style (fontId=2) distinct from the red one.
$ws.Range("C3").Font.ColorIndex = -4105

# D4 should pick up the existing red-font style (fontId=1), matching D3's look.
$ws.Range("D4").Font.Color = 255

# D7 gets a new numeric value.
$ws.Range("D7").Value = 34

# Update the active selection to E7.
$ws.Range("E7").Select()
